$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("FCREY")

# Row 8
$ws.Range("D8").Value = 6416600
$ws.Range("E8").Value = 6367800
$ws.Range("F8").Value = 6100200
$ws.Range("G8").Value = 5867800
$ws.Range("H8").Value = 5691700
$ws.Range("I8").Value = 5770300
$ws.Range("J8").Value = 5988400

# Row 9
$ws.Range("D9").Value = 5037200
$ws.Range("E9").Value = 4958600
$ws.Range("F9").Value = 4584600
$ws.Range("G9").Value = 4439700
$ws.Range("H9").Value = 4264200
$ws.Range("I9").Value = 4299400
$ws.Range("J9").Value = 8983700

# Row 10
$ws.Range("D10").Value = 1379400
$ws.Range("E10").Value = 1409200
$ws.Range("F10").Value = 1515600
$ws.Range("G10").Value = 1428200
$ws.Range("H10").Value = 1427500
$ws.Range("I10").Value = 1470900
$ws.Range("J10").Value = -2995200

# Row 14
$ws.Range("D14").Value = 214800
$ws.Range("E14").Value = 176200
$ws.Range("F14").Value = -11500
$ws.Range("G14").Value = 107000
$ws.Range("H14").Value = 33900
$ws.Range("I14").Value = 7500
$ws.Range("J14").Value = 107700

# Row 15
$ws.Range("H15").Value = 1400
$ws.Range("I15").Value = 1400

# Row 17
$ws.Range("D17").Value = 6496500
$ws.Range("E17").Value = 6182900
$ws.Range("F17").Value = 5613100
$ws.Range("G17").Value = 5527000
$ws.Range("H17").Value = 5290600
$ws.Range("I17").Value = 5384800
$ws.Range("J17").Value = 5715400

# Row 18
$ws.Range("D18").Value = -79900
$ws.Range("E18").Value = 185000
$ws.Range("F18").Value = 487100
$ws.Range("G18").Value = 340800
$ws.Range("H18").Value = 401100
$ws.Range("I18").Value = 385500
$ws.Range("J18").Value = 273000

# Row 20
$ws.Range("D20").Value = -19600
$ws.Range("G20").Value = -5400
$ws.Range("H20").Value = -6100
$ws.Range("I20").Value = -8100
$ws.Range("J20").Value = 88100

# Row 21
$ws.Range("D21").Value = 45800
$ws.Range("E21").Value = 318200
$ws.Range("F21").Value = 614200
$ws.Range("G21").Value = 472000
$ws.Range("H21").Value = 532900
$ws.Range("I21").Value = 526900
$ws.Range("J21").Value = 517400

# Row 22
$ws.Range("D22").Value = 86700
$ws.Range("E22").Value = 70500
$ws.Range("F22").Value = 73200
$ws.Range("G22").Value = 80600
$ws.Range("H22").Value = 82000
$ws.Range("I22").Value = 91500
$ws.Range("J22").Value = 191100

# Row 23
$ws.Range("D23").Value = -186300
$ws.Range("E23").Value = 109800
$ws.Range("F23").Value = 409200
$ws.Range("G23").Value = 254700
$ws.Range("H23").Value = 313000
$ws.Range("I23").Value = 285900
$ws.Range("J23").Value = 170100

# Row 24
$ws.Range("D24").Value = -65000
$ws.Range("E24").Value = 38600
$ws.Range("F24").Value = 88800
$ws.Range("G24").Value = 65000
$ws.Range("H24").Value = 75200
$ws.Range("I24").Value = 57600
$ws.Range("J24").Value = 39300

# Row 26
$ws.Range("D26").Value = -121300
$ws.Range("E26").Value = 71100
$ws.Range("F26").Value = 320500
$ws.Range("G26").Value = 189700
$ws.Range("H26").Value = 237800
$ws.Range("I26").Value = 228300
$ws.Range("J26").Value = 130800

# Row 27
$ws.Range("D27").Value = -128700
$ws.Range("E27").Value = 63700
$ws.Range("F27").Value = 313000
$ws.Range("G27").Value = 182900
$ws.Range("H27").Value = 229700
$ws.Range("I27").Value = 220900
$ws.Range("J27").Value = 125300

# Row 32
$ws.Range("D32").Value = 19600
$ws.Range("G32").Value = 5400
$ws.Range("H32").Value = 6100
$ws.Range("I32").Value = 8100
$ws.Range("J32").Value = -88100

# Row 33
$ws.Range("D33").Value = -128700
$ws.Range("E33").Value = 63700
$ws.Range("F33").Value = 313000
$ws.Range("G33").Value = 182900
$ws.Range("H33").Value = 229700
$ws.Range("I33").Value = 220900
$ws.Range("J33").Value = 125300

# Row 35
$ws.Range("D35").Value = -128700
$ws.Range("E35").Value = 63700
$ws.Range("F35").Value = 313000
$ws.Range("G35").Value = 182900
$ws.Range("H35").Value = 229700
$ws.Range("I35").Value = 220900
$ws.Range("J35").Value = 125300

# Row 41
$ws.Range("D41").Value = 613100
$ws.Range("E41").Value = 139600
$ws.Range("F41").Value = 99600
$ws.Range("G41").Value = 104300
$ws.Range("H41").Value = 51500
$ws.Range("I41").Value = 144300
$ws.Range("J41").Value = 174100

# Row 42
$ws.Range("D42").Value = 287900
$ws.Range("E42").Value = 8800
$ws.Range("F42").Value = 141600
$ws.Range("G42").Value = 50100
$ws.Range("H42").Value = 39300
$ws.Range("I42").Value = 22400
$ws.Range("J42").Value = 53500

# Row 43
$ws.Range("D43").Value = 2281800
$ws.Range("E43").Value = 1043300
$ws.Range("F43").Value = 924100
$ws.Range("G43").Value = 1037900
$ws.Range("H43").Value = 986400
$ws.Range("I43").Value = 932200
$ws.Range("J43").Value = 1949800

# Row 44
$ws.Range("D44").Value = 2368500
$ws.Range("E44").Value = 1119200
$ws.Range("F44").Value = 1002000
$ws.Range("G44").Value = 1020300
$ws.Range("H44").Value = 922800
$ws.Range("I44").Value = 916700
$ws.Range("J44").Value = 1943100

# Row 45
$ws.Range("D45").Value = 4100
$ws.Range("E45").Value = 5400
$ws.Range("F45").Value = 15600
$ws.Range("G45").Value = 4100
$ws.Range("H45").Value = 4100
$ws.Range("I45").Value = 10800
$ws.Range("J45").Value = 81300

# Row 46
$ws.Range("D46").Value = 2817000
$ws.Range("E46").Value = 2316400
$ws.Range("F46").Value = 2182900
$ws.Range("G46").Value = 2216800
$ws.Range("H46").Value = 2004000
$ws.Range("I46").Value = 1943100
$ws.Range("J46").Value = 2108400

# Row 47
$ws.Range("D47").Value = 101600
$ws.Range("E47").Value = 100300
$ws.Range("F47").Value = 92800
$ws.Range("G47").Value = 67800
$ws.Range("H47").Value = 90800
$ws.Range("I47").Value = 113800
$ws.Range("J47").Value = 302800

# Row 48
$ws.Range("D48").Value = 3036600
$ws.Range("E48").Value = 1500000
$ws.Range("F48").Value = 1343500
$ws.Range("G48").Value = 1505400
$ws.Range("H48").Value = 1418000
$ws.Range("I48").Value = 1758800
$ws.Range("J48").Value = 3181500

# Row 49
$ws.Range("D49").Value = 2284500
$ws.Range("E49").Value = 1136800
$ws.Range("F49").Value = 1154500
$ws.Range("G49").Value = 1151100
$ws.Range("H49").Value = 1103600
$ws.Range("I49").Value = 1521700
$ws.Range("J49").Value = 2387500

# Row 52
$ws.Range("D52").Value = 227000
$ws.Range("E52").Value = 145000
$ws.Range("F52").Value = 149700
$ws.Range("G52").Value = 140900
$ws.Range("H52").Value = 86000
$ws.Range("I52").Value = 87400
$ws.Range("J52").Value = 72500

# Row 54
$ws.Range("D54").Value = 5806200
$ws.Range("E54").Value = 5198500
$ws.Range("F54").Value = 4923400
$ws.Range("G54").Value = 5081900
$ws.Range("H54").Value = 4702500
$ws.Range("I54").Value = 4827200
$ws.Range("J54").Value = 5067000

# Row 57
$ws.Range("D57").Value = 1048100
$ws.Range("E57").Value = 952600
$ws.Range("F57").Value = 909200
$ws.Range("G57").Value = 890900
$ws.Range("H57").Value = 834000
$ws.Range("I57").Value = 827200
$ws.Range("J57").Value = 1692400

# Row 58
$ws.Range("D58").Value = 125300
$ws.Range("E58").Value = 182200
$ws.Range("F58").Value = 279800
$ws.Range("G58").Value = 230300
$ws.Range("H58").Value = 93500
$ws.Range("I58").Value = 102300
$ws.Range("J58").Value = 617900

# Row 59
$ws.Range("D59").Value = 956600
$ws.Range("E59").Value = 217500
$ws.Range("F59").Value = 164000
$ws.Range("G59").Value = 197800
$ws.Range("H59").Value = 153800
$ws.Range("I59").Value = 125300
$ws.Range("J59").Value = 298800

# Row 60
$ws.Range("D60").Value = 1705900
$ws.Range("E60").Value = 1352300
$ws.Range("F60").Value = 1353000
$ws.Range("G60").Value = 1319100
$ws.Range("H60").Value = 1081300
$ws.Range("I60").Value = 1054900
$ws.Range("J60").Value = 1311600

# Row 61
$ws.Range("D61").Value = 1187700
$ws.Range("E61").Value = 1289300
$ws.Range("F61").Value = 907200
$ws.Range("G61").Value = 1093500
$ws.Range("H61").Value = 1113800
$ws.Range("I61").Value = 1193100
$ws.Range("J61").Value = 1199900

# Row 62
$ws.Range("D62").Value = 106400
$ws.Range("E62").Value = 131400
$ws.Range("F62").Value = 147700
$ws.Range("G62").Value = 155800
$ws.Range("H62").Value = 167300
$ws.Range("I62").Value = 171400
$ws.Range("J62").Value = 412600

# Row 66
$ws.Range("D66").Value = 3016200
$ws.Range("E66").Value = 2789300
$ws.Range("F66").Value = 2422700
$ws.Range("G66").Value = 2586700
$ws.Range("H66").Value = 2386200
$ws.Range("I66").Value = 2443100
$ws.Range("J66").Value = 2750000

# Row 72
$ws.Range("D72").Value = 611800
$ws.Range("E72").Value = 832600
$ws.Range("F72").Value = 956600
$ws.Range("G72").Value = 822500
$ws.Range("H72").Value = 804900
$ws.Range("I72").Value = 731000
$ws.Range("J72").Value = 1335400

# Row 76
$ws.Range("D76").Value = 2789900
$ws.Range("E76").Value = 2409200
$ws.Range("F76").Value = 2500700
$ws.Range("G76").Value = 2495200
$ws.Range("H76").Value = 2316400
$ws.Range("I76").Value = 2384100
$ws.Range("J76").Value = 2317100

# Row 81
$ws.Range("D81").Value = -128700
$ws.Range("E81").Value = 63700
$ws.Range("F81").Value = 313000
$ws.Range("G81").Value = 182900
$ws.Range("H81").Value = 229700
$ws.Range("I81").Value = 220900
$ws.Range("J81").Value = 125300

# Row 83
$ws.Range("D83").Value = 145000
$ws.Range("E83").Value = 137500
$ws.Range("F83").Value = 131400
$ws.Range("G83").Value = 136200
$ws.Range("H83").Value = 137500
$ws.Range("I83").Value = 149100
$ws.Range("J83").Value = 155800

# Row 89
$ws.Range("D89").Value = 268300
$ws.Range("E89").Value = 164600
$ws.Range("F89").Value = 447200
$ws.Range("G89").Value = 389600
$ws.Range("H89").Value = 331300
$ws.Range("I89").Value = 378700
$ws.Range("J89").Value = 303500

# Row 91
$ws.Range("D91").Value = -206000
$ws.Range("E91").Value = -216100
$ws.Range("F91").Value = -203300
$ws.Range("G91").Value = -188300
$ws.Range("H91").Value = -176200
$ws.Range("I91").Value = -157900
$ws.Range("J91").Value = -176800

# Row 94
$ws.Range("D94").Value = -154500
$ws.Range("E94").Value = -411200
$ws.Range("F94").Value = -32500
$ws.Range("G94").Value = -145000
$ws.Range("H94").Value = -155100
$ws.Range("I94").Value = -105000
$ws.Range("J94").Value = -227600

# Row 96
$ws.Range("D96").Value = -83300
$ws.Range("E96").Value = -166700
$ws.Range("F96").Value = -159200
$ws.Range("G96").Value = -162600
$ws.Range("H96").Value = -151800
$ws.Range("I96").Value = -140900
$ws.Range("J96").Value = -136200

# Row 100
$ws.Range("D100").Value = 184300
$ws.Range("E100").Value = 155100
$ws.Range("F100").Value = -318400
$ws.Range("G100").Value = -187700
$ws.Range("H100").Value = -164600
$ws.Range("I100").Value = -302800
$ws.Range("J100").Value = -41300

# Row 101
$ws.Range("D101").Value = 4100
$ws.Range("E101").Value = -1400
$ws.Range("F101").Value = -9500
$ws.Range("G101").Value = 6800
$ws.Range("H101").Value = -4100
$ws.Range("I101").Value = -1400
$ws.Range("J101").Value = 1400

# Row 102
$ws.Range("D102").Value = 302200
$ws.Range("E102").Value = -92800
$ws.Range("F102").Value = 86700
$ws.Range("G102").Value = 63700
$ws.Range("H102").Value = 7500
$ws.Range("I102").Value = -30500
$ws.Range("J102").Value = 35900
